$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Ark1")
$ws2 = $wb.Worksheets.Item("Ark2")

# --- Time constraint: Ark1!D5 order's product-type window changes
# from "fisk3til4" (delivery in 3-4 days) to "fisk1til2" (delivery in 1-2 days)
$ws1.Range("D5").Value = "fisk1til2"

# --- Cross docking constraint: Ark2!C2 category for KristineKunden's
# out-of-country order switches from "b" to the new "contract" category
$ws2.Range("C2").Value = "contract"

# Reflect the active sheet/selection that results from this edit
$ws2.Range("C2").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D6").Select() | Out-Null
